$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") from 45192 to 45202 for all data rows (2 through 493)
for ($r = 2; $r -le 493; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# Row 493 gains an explicit row height (15, custom)
$ws.Rows.Item(493).RowHeight = 15

# Append a new data row 494
$ws.Cells.Item(494, 1).Value = "A 45700-2023"

$ws.Cells.Item(494, 2).Value = 45195
$ws.Cells.Item(494, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(494, 3).Value = 45202
$ws.Cells.Item(494, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(494, 4).Value = "ÖSTERGÖTLANDS LÄN"
$ws.Cells.Item(494, 5).Value = "ÅTVIDABERG"

$ws.Cells.Item(494, 7).Value = 2.6
$ws.Cells.Item(494, 8).Value = 0
$ws.Cells.Item(494, 9).Value = 0
$ws.Cells.Item(494, 10).Value = 0
$ws.Cells.Item(494, 11).Value = 0
$ws.Cells.Item(494, 12).Value = 0
$ws.Cells.Item(494, 13).Value = 0
$ws.Cells.Item(494, 14).Value = 0
$ws.Cells.Item(494, 15).Value = 0
$ws.Cells.Item(494, 16).Value = 0
$ws.Cells.Item(494, 17).Value = 0

$ws.Cells.Item(494, 18).WrapText = $true
